$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Round the coordinate values to integers
$ws.Range("Q2").Value = 596310
$ws.Range("R2").Value = 6572446

# Remove the "Starttid" (Z2) and "Sluttid" (AB2) time values entirely
$ws.Range("Z2").ClearContents()
$ws.Range("AB2").ClearContents()
